$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Sheet view: zoom / top-left cell / selection
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 140
$ws.Application.Goto($ws.Range("A5"))
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("B11").Select()

# ---------------------------------------------------------------------------
# 2. Text (shared-string) cell values.
#    Order matters only in that it controls the order new distinct strings
#    are appended to the workbook's shared-string table; cells that simply
#    reuse text already present elsewhere in the sheet can be written in any
#    order, so the "plain" ones are issued first and the genuinely new
#    pieces of text last, in the sequence they are first introduced.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Spring 2014"
$ws.Range("D3").Value = "Complete"
$ws.Range("A4").Value = "Get zip code-county mapping and write logic to lookup based on report from inpatient data"
$ws.Range("D4").Value = "X"
$ws.Range("D5").Value = "X"
$ws.Range("A6").Value = "Get bulk unemployment data with job distributions"
$ws.Range("A7").Value = "Explore visualizations using d3 or other libraries"
$ws.Range("B7").Value = "George"
$ws.Range("D7").Value = "X"
$ws.Range("A8").Value = "Analyze employment data with Pig"
$ws.Range("D8").Value = "X"
$ws.Range("A9").Value = "Explore CDC dataset for more health detail"
$ws.Range("A10").Value = "Develop paper draft based on the project proposal"
$ws.Range("B10").Value = "Chris"
$ws.Range("D10").Value = "X"
$ws.Range("B11").Value = "Chen"
$ws.Range("A12").Value = "Develop concrete results and commit to findings"
$ws.Range("B13").Value = "George"
$ws.Range("A14").Value = "Generate charts from reported analytic"
$ws.Range("B14").Value = "George"
$ws.Range("B15").Value = "Chris"
$ws.Range("B18").Value = "George"

$ws.Range("B8").Value = "Khen, Chris"
$ws.Range("D6").Value = "N/A"
$ws.Range("D9").Value = "N/A"
$ws.Range("A13").Value = "Choose functional attributes for regression (Weka/STATA)"
$ws.Range("A15").Value = "Update paper with draft of findings (RESULTS)"
$ws.Range("A16").Value = "Update paper with draft of future work"
$ws.Range("A17").Value = "Update paper with draft of Conclusions"
$ws.Range("A19").Value = "Presentation "
$ws.Range("A18").Value = "Presentation Slides and draft script"
$ws.Range("B6").Value = "Chris, Khen, George"
$ws.Range("B19").Value = "Chris, Khen, George"
$ws.Range("B4").Value = "Khen"
$ws.Range("B12").Value = "!"
$ws.Range("A11").Value = "Breakdown existing census employment data by fields and consolidate into combined results"

# ---------------------------------------------------------------------------
# 3. Date values in column C (numeric, formatted via existing style "3")
# ---------------------------------------------------------------------------
$ws.Range("C6").Value = 41739
$ws.Range("C7").Value = 41746
$ws.Range("C8").Value = 41746
$ws.Range("C9").Value = 41739
$ws.Range("C10").Value = 41736
$ws.Range("C11").Value = 41751
$ws.Range("C12").Value = 41752
$ws.Range("C13").Value = 41753
$ws.Range("C14").Value = 41757
$ws.Range("C15").Value = 41760
$ws.Range("C16").Value = 41760
$ws.Range("C17").Value = 41760
$ws.Range("C18").Value = 41760
$ws.Range("C19").Value = 41767

# ---------------------------------------------------------------------------
# 4. Apply the existing date format (style id 3, taken from C4) to the newly
#    populated C14:C19 cells, and the existing centered "flag" format
#    (style id 4, taken from D4) to the newly populated D6:D10 cells - this
#    reuses the pre-existing style slots instead of fabricating new ones.
# ---------------------------------------------------------------------------
$ws.Range("C4").Copy()
$ws.Range("C14:C19").PasteSpecial(-4122)

$ws.Range("D4").Copy()
$ws.Range("D6:D10").PasteSpecial(-4122)

$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5. Remove the old D12 flag cell entirely (row 12 no longer has a column D
#    entry in the revised sheet).
# ---------------------------------------------------------------------------
$ws.Range("D12").Clear()

# ---------------------------------------------------------------------------
# 6. Highlight A12 with a new yellow fill + wrapped text style.
# ---------------------------------------------------------------------------
$ws.Range("A12").Interior.Color = 65535
$ws.Range("A12").WrapText = $true

# ---------------------------------------------------------------------------
# 7. Row-height tweaks that differ from the sheet's previous layout.
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 45
$ws.Rows.Item(15).RowHeight = 30

Write-Output "done"
